# Replace the run(s) covered by $rng with the raw WordprocessingML in
# $runsXml, wrapped in a minimal Flat-OPC "pkg:package" payload so
# Range.InsertXML can parse it. Using a sub-range that matches exactly
# the span of the run(s) being replaced (instead of the whole paragraph)
# leaves any sibling runs -- in particular the empty <w:r/> marker runs
# this document uses -- untouched.
function Set-RunsXml($rng, [string]$runsXml) {
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml) | Out-Null
}

# Same idea, but replaces a span of whole paragraphs (<w:p>...</w:p>) with
# a different set of paragraphs -- used where the number of list items
# changes.
function Set-ParasXml($rng, [string]$parasXml) {
    $xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
        $parasXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($xml) | Out-Null
}

$d = $word.ActiveDocument

# 1) Title (Heading1 at the very top of the document) -- this paragraph
#    has no leading empty run, so a plain Find/Replace is safe.
$d.Content.Find.Execute(
    "Play Mythological Mayhem Supreme Streaks Free - Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Mythological Mayhem Supreme Streaks for Free", 2) | Out-Null

# 2) "What we like" bullet list. The 4 existing items become 4 different
#    items (not a 1:1 rename), so replace the whole paragraph span.
$pLikeFirst = $d.Paragraphs.Item(35)   # "Innovative graphics and animations"
$pLikeLast  = $d.Paragraphs.Item(38)   # "117,649 ways to win"
$rngLike = $d.Range($pLikeFirst.Range.Start, $pLikeLast.Range.End)
$pPrBullet = '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>'
$likeItems = @(
    "Engaging mythological theme",
    "Innovative gameplay features",
    "Impressive graphics and animations",
    "Exciting bonus stages"
)
$likeParasXml = ""
foreach ($it in $likeItems) {
    $likeParasXml += "<w:p>$pPrBullet<w:r/><w:r><w:t>$it</w:t></w:r></w:p>"
}
Set-ParasXml $rngLike $likeParasXml

# 3) "What we don't like" bullet list -- simple text swaps, same count.
$pRtp = $d.Paragraphs.Item(40)         # "RTP is slightly below industry average"
$rngRtp = $d.Range($pRtp.Range.Start, $pRtp.Range.End)
Set-RunsXml $rngRtp '<w:r><w:t>RTP slightly below average</w:t></w:r>'

$pLim = $d.Paragraphs.Item(41)         # "Limited range of symbols"
$rngLim = $d.Range($pLim.Range.Start, $pLim.Range.End)
Set-RunsXml $rngLim '<w:r><w:t>Limited number of worlds</w:t></w:r>'

# 4) Bold "title" paragraph near the bottom of the doc.
$pBold = $d.Paragraphs.Item(42)
$rngBold = $d.Range($pBold.Range.Start, $pBold.Range.End)
Set-RunsXml $rngBold '<w:r><w:rPr><w:b/></w:rPr><w:t>Play Mythological Mayhem Supreme Streaks for Free</w:t></w:r>'

# 5) Italic meta-description paragraph (last paragraph in the document).
$pItalic = $d.Paragraphs.Item(43)
$rngItalic = $d.Range($pItalic.Range.Start, $pItalic.Range.End)
Set-RunsXml $rngItalic '<w:r><w:rPr><w:i/></w:rPr><w:t>Experience ancient Greece with this thrilling online slot game. Play for free and win big!</w:t></w:r>'
